$wb = $excel.ActiveWorkbook

# The workbook has two sheets that carry the same "漫展" listing data:
#   "展览"     (Exhibitions)
#   "全部类型" (All types - mirrors the exhibitions sheet)
# Both need the same numeric updates to column F ("想去人数" / interest count).
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 513
    $ws.Range("F6").Value = 84
    $ws.Range("F7").Value = 712
}
